$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and date range) ---
$ws.Range("A8").Value = "Volume 31   Number  23"
$ws.Range("C9").Value = "Report Covering the Week  6/3/2024  Through  6/9/2024"

# --- Crime Complaints data table (rows 14-31) ---
# Cells whose style/type changes (text <-> number) - copy a reference cell's format first
$ws.Range("I14").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("D15").Value = 1
$ws.Range("L14").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("E15").Value = -100
$ws.Range("I14").Copy()
$ws.Range("G15").PasteSpecial(-4122)
$ws.Range("G15").Value = 1
$ws.Range("L14").Copy()
$ws.Range("H15").PasteSpecial(-4122)
$ws.Range("H15").Value = -100
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "0"
$ws.Range("F14").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("I14").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("C20").Value = 1
$ws.Range("I14").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("D22").Value = 1
$ws.Range("L14").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("E22").Value = -100
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "0"
$ws.Range("F14").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("I14").Copy()
$ws.Range("C25").PasteSpecial(-4122)
$ws.Range("C25").Value = 1
$ws.Range("I14").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("C27").Value = 3
$ws.Range("I14").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("D27").Value = 2
$ws.Range("L14").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("E27").Value = 50
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0"
$ws.Range("F14").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "***.*"
$ws.Range("F14").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "0"
$ws.Range("F14").Copy()
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "0"
$ws.Range("F14").Copy()
$ws.Range("C30").PasteSpecial(-4122)
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = "0"
$ws.Range("F14").Copy()
$ws.Range("C31").PasteSpecial(-4122)

# Cells whose value changes only (style/type unchanged)
$ws.Range("N14").Value = -75
$ws.Range("J15").Value = 7
$ws.Range("K15").Value = 42.857142857142
$ws.Range("N15").Value = -28.571428571428
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 200
$ws.Range("F16").Value = 8
$ws.Range("G16").Value = 6
$ws.Range("H16").Value = 33.333333333333
$ws.Range("I16").Value = 47
$ws.Range("J16").Value = 52
$ws.Range("K16").Value = -9.615384615384
$ws.Range("L16").Value = 34.285714285714
$ws.Range("M16").Value = -14.545454545454
$ws.Range("N16").Value = -82.264150943396
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = -16.666666666666
$ws.Range("F17").Value = 20
$ws.Range("H17").Value = -13.043478260869
$ws.Range("I17").Value = 141
$ws.Range("J17").Value = 101
$ws.Range("K17").Value = 39.603960396039
$ws.Range("L17").Value = 39.603960396039
$ws.Range("M17").Value = 101.428571428571
$ws.Range("N17").Value = -24.598930481283
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = -100
$ws.Range("F18").Value = 4
$ws.Range("H18").Value = -33.333333333333
$ws.Range("J18").Value = 33
$ws.Range("K18").Value = -6.060606060606
$ws.Range("L18").Value = 6.896551724137
$ws.Range("M18").Value = -39.215686274509
$ws.Range("N18").Value = -90.127388535031
$ws.Range("C19").Value = 2
$ws.Range("D19").Value = 3
$ws.Range("E19").Value = -33.333333333333
$ws.Range("F19").Value = 14
$ws.Range("G19").Value = 14
$ws.Range("I19").Value = 70
$ws.Range("J19").Value = 78
$ws.Range("K19").Value = -10.256410256410
$ws.Range("L19").Value = -18.604651162790
$ws.Range("M19").Value = 75
$ws.Range("N19").Value = -40.677966101694
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 0
$ws.Range("G20").Value = 8
$ws.Range("H20").Value = -62.5
$ws.Range("J20").Value = 33
$ws.Range("K20").Value = -18.181818181818
$ws.Range("L20").Value = 28.571428571428
$ws.Range("N20").Value = -88.461538461538
$ws.Range("C21").Value = 11
$ws.Range("D21").Value = 13
$ws.Range("E21").Value = -15.384615384615
$ws.Range("G21").Value = 58
$ws.Range("H21").Value = -15.517241379310
$ws.Range("I21").Value = 328
$ws.Range("J21").Value = 304
$ws.Range("K21").Value = 7.894736842105
$ws.Range("L21").Value = 15.492957746478
$ws.Range("M21").Value = 24.714828897338
$ws.Range("N21").Value = -71.228070175438
$ws.Range("F22").Value = 1
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 8
$ws.Range("K22").Value = -25
$ws.Range("E23").Value = -100
$ws.Range("F23").Value = 9
$ws.Range("G23").Value = 7
$ws.Range("H23").Value = 28.571428571428
$ws.Range("J23").Value = 40
$ws.Range("K23").Value = -2.5
$ws.Range("L23").Value = 18.181818181818
$ws.Range("M23").Value = 116.666666666667
$ws.Range("C24").Value = 14
$ws.Range("D24").Value = 13
$ws.Range("E24").Value = 7.692307692307
$ws.Range("F24").Value = 50
$ws.Range("G24").Value = 53
$ws.Range("H24").Value = -5.660377358490
$ws.Range("I24").Value = 264
$ws.Range("J24").Value = 250
$ws.Range("K24").Value = 5.6
$ws.Range("L24").Value = 8.641975308641
$ws.Range("M24").Value = 68.152866242038
$ws.Range("E25").Value = -50
$ws.Range("G25").Value = 12
$ws.Range("H25").Value = -66.666666666666
$ws.Range("I25").Value = 21
$ws.Range("J25").Value = 47
$ws.Range("K25").Value = -55.319148936170
$ws.Range("L25").Value = -34.375
$ws.Range("C26").Value = 9
$ws.Range("D26").Value = 10
$ws.Range("E26").Value = -10
$ws.Range("F26").Value = 29
$ws.Range("G26").Value = 40
$ws.Range("H26").Value = -27.5
$ws.Range("I26").Value = 165
$ws.Range("J26").Value = 171
$ws.Range("K26").Value = -3.508771929824
$ws.Range("L26").Value = 3.773584905660
$ws.Range("M26").Value = -6.779661016949
$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = -25
$ws.Range("I27").Value = 17
$ws.Range("J27").Value = 18
$ws.Range("K27").Value = -5.555555555555
$ws.Range("L27").Value = -5.555555555555
$ws.Range("F28").Value = 6
$ws.Range("H28").Value = 20
$ws.Range("I28").Value = 20
$ws.Range("K28").Value = 5.263157894736
$ws.Range("L28").Value = 0
$ws.Range("F31").Value = 1

$ws.Application.CutCopyMode = $false

